# Row-level price/volume refresh, plus a coin re-ranking shift in rows 6-18
# (GateToken moved to the top of that block; FTXToken, MXToken, BTSEToken, etc.
# each shifted down one row and LEO/GateToken rotated to the bottom).
$updates = @(
    @{ Row=2; D='314.56'; E='1.40%' }
    @{ Row=3; D='40.90'; E='-1.05%' }
    @{ Row=4; D='5.131'; E='0.27%' }
    @{ Row=5; D='0.07634'; E='-0.85%' }
    @{ Row=6; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='4.336'; E='0.65%' }
    @{ Row=7; B='FTXToken'; C='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D='1.685'; E='3.63%' }
    @{ Row=8; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='0.9336'; E='1.11%' }
    @{ Row=9; B='BTSEToken'; C='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; D='2.425'; E='-1.74%' }
    @{ Row=10; B='LiechtensteinCryptoassetsExchange'; C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D='0.1247'; E='2.74%' }
    @{ Row=11; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1818'; E='-1.08%' }
    @{ Row=12; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.09021'; E='-1.81%' }
    @{ Row=13; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.04141'; E='-1.70%' }
    @{ Row=14; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.1056'; E='0.66%' }
    @{ Row=15; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001295'; E='3.05%' }
    @{ Row=16; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.005809'; E='-0.16%' }
    @{ Row=17; B='UpBots'; C='https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'; D='0.007491'; E='1,897.31%' }
    @{ Row=18; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.355'; E='0.15%' }
    @{ Row=19; D='0.3358'; E='1.63%' }
    @{ Row=20; D='8.386'; E='21.34%' }
    @{ Row=21; D='0.1347'; E='-2.95%' }
    @{ Row=22; D='0.2742'; E='2.49%' }
    @{ Row=23; D='0.04042'; E='-0.14%' }
    @{ Row=24; D='0.001267'; E='0.50%' }
    @{ Row=25; D='0.004053'; E='-1.25%' }
    @{ Row=26; D='0.0001277'; E='0.72%' }
    @{ Row=38; D='0.02483'; E='0.42%' }
    @{ Row=39; D='0.05193'; E='-1.32%' }
    @{ Row=40; D='0.007807'; E='-0.23%' }
    @{ Row=41; E='-1.17%' }
    @{ Row=42; D='0.007383'; E='8.74%' }
    @{ Row=43; D='0.002159'; E='16.25%' }
    @{ Row=44; D='0.008249'; E='0.74%' }
    @{ Row=45; D='0.3144'; E='1.38%' }
    @{ Row=46; D='0.00006671'; E='-0.72%' }
    @{ Row=47; D='0.00000000754'; E='0.76%' }
    @{ Row=48; D='0.2503'; E='21.83%' }
    @{ Row=49; D='0.004223'; E='3.22%' }
    @{ Row=50; D='0.00002111'; E='0.76%' }
    @{ Row=51; D='0.0002011'; E='0.76%' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }

    # Columns D (Price) and E (Volume(1h)) hold numeric-looking text
    # ("314.56", "-1.05%", "1,897.31%", ...). A bare .Value assignment would
    # let Excel auto-convert these into real numbers/percentages (dropping
    # significant trailing zeros and reformatting). Prefixing with a single
    # quote forces a text entry matching the original inlineStr cells, and
    # re-applying the Normal style clears the quote-prefix formatting mark
    # Excel adds so the cell format stays identical to the untouched cells.
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey('E')) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.Value = "'" + $u.E
        $cell.Style = "Normal"
    }
}
